$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 185, shifting existing rows 185-196 down to 186-197.
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with the new weekly price record.
$ws.Cells.Item(185, 1).Value  = 3
$ws.Cells.Item(185, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(185, 3).Value  = "Coquimbo"
$ws.Cells.Item(185, 4).Value  = 44610
$ws.Cells.Item(185, 5).Value  = 5
$ws.Cells.Item(185, 6).Value  = 100112010
$ws.Cells.Item(185, 7).Value  = "Achicoria"
$ws.Cells.Item(185, 8).Value  = "Sin especificar"
$ws.Cells.Item(185, 9).Value  = "Primera"
$ws.Cells.Item(185, 10).Value = 48
$ws.Cells.Item(185, 11).Value = 7000
$ws.Cells.Item(185, 12).Value = 7000
$ws.Cells.Item(185, 13).Value = 7000
$ws.Cells.Item(185, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(185, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(185, 16).Value = 438
$ws.Cells.Item(185, 17).Value = 16
$ws.Cells.Item(185, 18).Value = "Hortaliza"
